$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.521.14"
$ws.Range("E2").Value = "  +0.55%  "

# Row 3
$ws.Range("D3").Value = "3.595.59"
$ws.Range("E3").Value = "  +0.78%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "'608.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "

# Row 6
$ws.Range("D6").Value = "'149.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.66%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.489"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.16%  "

# Row 9
$ws.Range("E9").Value = "  -0.46%  "

# Row 10
$ws.Range("D10").Value = "'8.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.14%  "

# Row 11
$ws.Range("E11").Value = "  +0.70%  "

# Row 12
$ws.Range("D12").Value = "4.204.20"
$ws.Range("E12").Value = "  +0.78%  "

# Row 13
$ws.Range("E13").Value = "  +0.13%  "

# Row 14
$ws.Range("E14").Value = "  -0.88%  "

# Row 15
$ws.Range("D15").Value = "3.594.44"
$ws.Range("E15").Value = "  +0.75%  "

# Row 16
$ws.Range("D16").Value = "66.573.01"
$ws.Range("E16").Value = "  +0.48%  "

# Row 17
$ws.Range("E17").Value = "  +0.94%  "

# Row 18
$ws.Range("E18").Value = "  +1.31%  "

# Row 19
$ws.Range("E19").Value = "  +2.49%  "

# Row 20
$ws.Range("E20").Value = "  +1.62%  "

# Row 21
$ws.Range("D21").Value = "'427.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.69%  "

# Row 22
$ws.Range("E22").Value = "  +1.00%  "

# Row 23
$ws.Range("D23").Value = "'78.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.56%  "

# Row 24
$ws.Range("E24").Value = "  -0.01%  "

# Row 25
$ws.Range("E25").Value = "  +1.37%  "

# Row 26
$ws.Range("D26").Value = "'8.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.29%  "

# Row 27
$ws.Range("D27").Value = "'9.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.98%  "

# Row 28
$ws.Range("E28").Value = "  +0.39%  "

# Row 29
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.19%  "

# Row 30
$ws.Range("D30").Value = "3.590.67"
$ws.Range("E30").Value = "  +0.80%  "

# Row 31
$ws.Range("E31").Value = "  +0.07%  "

# Row 32
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "'0.157"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.89%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'25.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.99%  "

# Row 34
$ws.Range("D34").Value = "'7.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.33%  "

# Row 36
$ws.Range("E36").Value = "  +0.32%  "

# Row 37
$ws.Range("E37").Value = "  -3.02%  "

# Row 38
$ws.Range("D38").Value = "'177.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.82%  "

# Row 39
$ws.Range("D39").Value = "'0.0856"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.44%  "

# Row 40
$ws.Range("E40").Value = "  +0.36%  "

# Row 41
$ws.Range("D41").Value = "'0.897"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "

# Row 42
$ws.Range("E42").Value = "  -2.74%  "

# Row 43
$ws.Range("E43").Value = "  +8.41%  "

# Row 44
$ws.Range("D44").Value = "'0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "

# Row 45
$ws.Range("D45").Value = "'25.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.27%  "

# Row 46
$ws.Range("E46").Value = "  -2.49%  "

# Row 47
$ws.Range("D47").Value = "'24.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.57%  "

# Row 48
$ws.Range("D48").Value = "'7.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.99%  "

# Row 49
$ws.Range("E49").Value = "  +0.91%  "

# Row 50
$ws.Range("D50").Value = "2.430.88"
$ws.Range("E50").Value = "  +5.14%  "

# Row 51
$ws.Range("D51").Value = "'0.235"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.23%  "
